$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9187717437744141
$ws.Range("B1").Value = 1.553135395050049
$ws.Range("C1").Value = 4.369341373443604
$ws.Range("D1").Value = 2.625953674316406
$ws.Range("E1").Value = 1.470887899398804
